# Updates the Malta Premier League 2023-2024 results sheet:
#  - Rows 19/20, 40/41, 48/49 and 75/76 had their match details (columns
#    F:V) swapped with each other (the "Indice"/A column and the
#    pais/torneio/temporada/data_partida columns A:E stay put - only the
#    match content moved rows).
#  - Two new fixtures (rows 78 and 79) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($rowA, $rowB) {
    $rangeA = $ws.Range("F$rowA`:V$rowA")
    $rangeB = $ws.Range("F$rowB`:V$rowB")

    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()

    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-MatchRows 19 20
Swap-MatchRows 40 41
Swap-MatchRows 48 49
Swap-MatchRows 75 76

# Append the two new fixtures below the existing last row (77), copying
# row 77's formatting (bold/bordered index column, date-time formatted
# match-date column) onto the new rows.
$ws.Range("A77:V77").Copy($ws.Range("A78:V78"))
$ws.Range("A77:V77").Copy($ws.Range("A79:V79"))

$row78 = @{
    A = 77; B = "malta"; C = "premier-league"; D = "2023-2024"; E = 45297.58333333334
    F = "Hamrun"; G = 1; H = "Floriana"; I = 0
    J = 2.1;  K = "05/01/2024 02:12"; L = 2.22; M = "06/01/2024 13:56"
    N = 2.9;  O = "05/01/2024 02:12"; P = 2.63; Q = "06/01/2024 13:56"
    R = 3.51; S = "05/01/2024 02:12"; T = 4.15; U = "06/01/2024 13:57"
    V = "https://www.betexplorer.com/football/malta/premier-league/hamrun-floriana/2mP09kzD/"
}

$row79 = @{
    A = 78; B = "malta"; C = "premier-league"; D = "2023-2024"; E = 45297.70833333334
    F = "Gzira"; G = 1; H = "Santa Lucia"; I = 2
    J = 1.44; K = "05/01/2024 05:12"; L = 1.71; M = "06/01/2024 16:58"
    N = 4.13; O = "05/01/2024 05:12"; P = 3.18; Q = "06/01/2024 16:59"
    R = 5.9;  S = "05/01/2024 05:12"; T = 5.85; U = "06/01/2024 16:58"
    V = "https://www.betexplorer.com/football/malta/premier-league/gzira-santa-lucia/tvO489kJ/"
}

foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")) {
    $ws.Range("$col`78").Value = $row78[$col]
    $ws.Range("$col`79").Value = $row79[$col]
}
